$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of annotations (id, data, questao_id, disciplina, assunto, anotacao)
$rows = @(
    @(22, "29/12/2025 04:28", 545,  "Conhecimentos Específicos", "Gestão da Manutenção e Confiabilidade", "Criar um formulários/flashcards para as fórmulas desse assunto"),
    @(23, "29/12/2025 04:45", 1007, "Estatística", "Medidas de Variabilidade", "Revisar medidas de variabilidade em tabelas de frequencias"),
    @(24, "29/12/2025 04:46", 1191, "Estatística", "Testes de Hipóteses", "p-value, faço nem ideia de como se calcula"),
    @(25, "29/12/2025 04:49", 1241, "Estatística", "ANOVA", "Faço nem ideia"),
    @(26, "29/12/2025 04:58", "954", "Inglês", "Pronouns", "Organizar esse texto")
)

$startRow = 13
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]

    if ($r -eq 17) {
        # questao_id stored as text for this row, matching source data
        $ws.Cells.Item($r, 3).NumberFormat = "@"
        $ws.Cells.Item($r, 3).Value = $data[2]
    } else {
        $ws.Cells.Item($r, 3).Value = $data[2]
    }

    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
}
